# ============================================================================
# Applies the "additional scraping" update:
#   1. Adds a new "Player Info" sheet (before "ODI Batting")
#   2. Renames MATCH_CARD_LINK -> MATCH_CODE in "ODI Batting" and "ODI Bowling"
#      and turns the full scorecard URL into just the bare match code
#   3. Removes the stray empty B-cells on the "did not bat" rows in
#      "ODI Batting"
#   4. Adds a new "ODI Batting Extra" sheet (after "ODI Bowling")
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: apply the same "header" look (bold, bordered, centered) that the
# existing sheets use for row 1.
# ----------------------------------------------------------------------------
function Format-HeaderRange($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# ----------------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to stay text (so that
# purely-numeric looking strings such as match codes or player ids are not
# silently turned into numbers by Excel's auto-detection).
# ----------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ============================================================================
# 1. Reorder / create the sheets
# ============================================================================
$oldBatting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($oldBatting)
$playerInfo.Name = "Player Info"

$bowling = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowling)
$battingExtra.Name = "ODI Batting Extra"

$batting = $wb.Worksheets.Item("ODI Batting")

# ============================================================================
# 2. Player Info sheet
# ============================================================================
$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
}
Format-HeaderRange($playerInfo.Range("A1:D1"))

Set-TextValue $playerInfo.Cells.Item(2, 1) "3845"
$playerInfo.Cells.Item(2, 2).Value = "Aaron James Finch"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# ============================================================================
# 3. ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE (column D)
# ============================================================================
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingLastRow = $batting.UsedRange.Rows.Count
$battingCol = $batting.Range("D2:D" + $battingLastRow)
$battingCol.NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $txt = $cell.Text
    if ($txt -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}
$battingCol.Style = "Normal"

# Remove the stray empty INNING_NUMBER cells on the "did not bat" rows
foreach ($r in 11, 41, 63, 70) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# ============================================================================
# 4. ODI Bowling sheet: MATCH_CARD_LINK -> MATCH_CODE (column B)
# ============================================================================
$bowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingLastRow = $bowling.UsedRange.Rows.Count
$bowlingCol = $bowling.Range("B2:B" + $bowlingLastRow)
$bowlingCol.NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowling.Cells.Item($r, 2)
    $txt = $cell.Text
    if ($txt -match "MatchCode=(\d+)") {
        $cell.Value = $matches[1]
    }
}
$bowlingCol.Style = "Normal"

# ============================================================================
# 5. ODI Batting Extra sheet
# ============================================================================
$beHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $beHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c)
    $cell.Value = $beHeaders[$c - 1]
}
Format-HeaderRange($battingExtra.Range("A1:F1"))

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$beData = @(
    @("4429", 2, "3", "0", "5.44%", "NO"),
    @("4430", 2, "8", "1", "35.27%", "NO"),
    @("4431", 2, "2", "0", "3.93%", "NO"),
    @("4435", $null, $null, $null, $null, "NO"),
    @("4436", $null, $null, $null, $null, "NO"),
    @("4437", 2, "7", "3", "25.95%", "NO"),
    @("4564", 2, "0", "1", "7.35%", "NO"),
    @("4565", $null, $null, $null, $null, "NO"),
    @("4567", 2, "0", "0", $null, "NO"),
    @("4594", $null, $null, $null, $null, "NO"),
    @("4597", 2, "2", "0", "7.41%", "NO"),
    @("4600", $null, $null, $null, $null, "NO"),
    @("4601", 2, "0", "0", $null, "NO"),
    @("4603", 2, "0", "0", $null, "NO"),
    @("4644", 1, "2", "0", "7.46%", "NO"),
    @("4645", 2, "0", "0", "1.00%", "NO"),
    @("4646", 2, "1", "0", "3.55%", "NO"),
    @("4647", $null, $null, $null, $null, "NO"),
    @("4648", 2, "0", "0", $null, "NO"),
    @("4649", 2, "0", "0", "1.87%", "NO")
)

for ($i = 0; $i -lt $beData.Length; $i++) {
    $row = $i + 2
    $values = $beData[$i]

    Set-TextValue $battingExtra.Cells.Item($row, 1) $values[0]

    $posCell = $battingExtra.Cells.Item($row, 2)
    if ($null -ne $values[1]) {
        $posCell.Value = $values[1]
    } else {
        $posCell.ClearContents()
        $posCell.Value = ""
    }

    foreach ($colInfo in @(@(3, $values[2]), @(4, $values[3]), @(5, $values[4]))) {
        $col = $colInfo[0]
        $val = $colInfo[1]
        $cell = $battingExtra.Cells.Item($row, $col)
        if ($null -ne $val) {
            Set-TextValue $cell $val
        } else {
            $cell.ClearContents()
            $cell.Value = ""
        }
    }

    $battingExtra.Cells.Item($row, 6).Value = $values[5]
}
